$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.358.75'

$ws.Range('D3').Value = '2.278.62'

$ws.Range('E3').Value = '  +0.58%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.02'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').Value = '  +1.18%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.62'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').Value = '  +5.96%  '

$ws.Range('E9').Value = '  +2.21%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.55'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').Value = '  +9.69%  '

$ws.Range('E11').Value = '  +0.07%  '

$ws.Range('E12').Value = '  -1.17%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.69'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').Value = '  +0.66%  '

$ws.Range('D14').Value = '2.628.85'

$ws.Range('E14').Value = '  +0.48%  '

$ws.Range('E15').Value = '  +1.37%  '

$ws.Range('D16').Value = '2.268.35'

$ws.Range('E16').Value = '  -0.92%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.796'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').Value = '  +3.48%  '

$ws.Range('D18').Value = '42.240.56'

$ws.Range('E18').Value = '  +1.41%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.56'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').Value = '  +0.58%  '

$ws.Range('E20').Value = '  +0.56%  '

$ws.Range('E21').Value = '  +0.51%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.60'
$ws.Range('D22').Style = 'Normal'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.51'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').Value = '  +0.34%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.59'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').Value = '  +0.29%  '

$ws.Range('E25').Value = '  +1.14%  '

$ws.Range('E26').Value = '  +0.02%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.83'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E28').Value = '  +6.01%  '

$ws.Range('E29').Value = '  -0.51%  '

$ws.Range('E30').Value = '  +1.24%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '159.85'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').Value = '  -0.15%  '

$ws.Range('E32').Value = '  -0.16%  '

$ws.Range('E33').Value = '  -0.07%  '

$ws.Range('E34').Value = '  +3.70%  '

$ws.Range('E35').Value = '  -0.37%  '

$ws.Range('E36').Value = '  +0.17%  '

$ws.Range('E37').Value = '  +0.79%  '

$ws.Range('E39').Value = '  +3.20%  '

$ws.Range('E41').Value = '  +4.95%  '

$ws.Range('E42').Value = '  +14.16%  '

$ws.Range('D43').Value = '1.999.22'

$ws.Range('E43').Value = '  -0.67%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0285'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').Value = '  +2.33%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.86'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').Value = '  -2.00%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.98'
$ws.Range('D46').Style = 'Normal'

$ws.Range('E46').Value = '  +3.46%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.00'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').Value = '  -3.63%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '52.97'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').Value = '  +1.46%  '

$ws.Range('E49').Value = '  +0.71%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.08'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E50').Value = '  +0.19%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '91.87'
$ws.Range('D51').Style = 'Normal'

$ws.Range('E51').Value = '  +0.84%  '
